# Insert a new weekly record above the current row 21 ("Vega Central
# Mapocho de Santiago - Arveja Verde"), shifting all existing rows 21-48
# down to 22-49 and growing the used range to A1:R49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 21 (pushes 21..48 down to 22..49).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new data point.
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44413
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112022
$ws.Range("G21").Value = "Arveja Verde"
$ws.Range("H21").Value = "Perfection"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 16
$ws.Range("K21").Value = 33000
$ws.Range("L21").Value = 34000
$ws.Range("M21").Value = 33500
$ws.Range("N21").Value = "$/malla 25 kilos"
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 1340
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
